$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.027798652648926
$ws.Range("B1").Value = 1.907064914703369
$ws.Range("C1").Value = 7.789341449737549
$ws.Range("D1").Value = 2.075407028198242
$ws.Range("E1").Value = 0.552586555480957
